$wb = $excel.ActiveWorkbook

# Sheet "A5" -> corresponds to the evidence row for Task A5
$wsA5 = $wb.Worksheets.Item("A5")
$wsA5.Range("A2").Value = "DDB1E78D79A3F02BD7262F37418B353994550D194D93420B61B1FB8B592F1287"
$wsA5.Range("B2").Value = "wasm.juno1aff4yjmkp9lqce596pujxnknfthrpyj3c0n7dmm7q6rr0gz37pzsz49e9z"
$wsA5.Range("C2").Value = "landerosGoN01"
$wsA5.Range("D2").Value = "uni-6"
$wsA5.Range("B8").Select() | Out-Null

# Sheet "A6" -> corresponds to the evidence row for Task A6
$wsA6 = $wb.Worksheets.Item("A6")
$wsA6.Range("A2").Value = "DB0ACAFD86DDED8DAB7D92315A6FF010B041E051CC49DFA950D4F6576EAB39C4"
$wsA6.Range("B2").Value = "ibc/E8497616125E97550485C786E64AD14DE49AFBE17BD7F5B1933804A88D9DD7C9"
$wsA6.Range("C2").Value = "landerosGoN01"
$wsA6.Range("D2").Value = "uptick_7000-2"
$wsA6.Rows.Item(1).AutoFit() | Out-Null
$wsA6.Rows.Item(2).AutoFit() | Out-Null
$wsA6.Range("C8").Select() | Out-Null
